## i230--EIT2 ta and mpar suite update
## Applies the content changes described in the commit: remove the
## --devkit=iCE40UP5K-CM225I switch from the TA "run-map-trce" command
## line on the suite sheet (it now lives on the case sheet as its own
## per-row MPAR option column), populate the new MPAR "cmd" column (M)
## and the new "mpar" marker column (Q) for the mpar case rows, and add
## three new "hold correction" case rows.

$wb = $excel.ActiveWorkbook
$suite = $wb.Worksheets.Item("suite")
$case = $wb.Worksheets.Item("case")

# New MPAR_option cmd column for every existing mpar-suite case row
$case.Range("M3:M24").Value = "cmd = --devkit=iCE40UP5K-CM225I"

# "suite" sheet: drop the CM225I devkit switch from the run-map-trce cmd
# (the switch now lives per-row in the case sheet's new M column above)
$suite.Range("B7").Value = "cmd =  python DEV/bin/run_diamondng.py --synthesis=lse --run-map-trce --run-par-trce --run-export-bitstream"

# Three new "hold correction" case rows
$case.Range("D26").Value = "hold correction"
$case.Range("E26").Value = "hold_correction/hold_correction1"

$case.Range("D27").Value = "hold correction"
$case.Range("E27").Value = "hold_correction/hold_correction2"

$case.Range("D28").Value = "hold correction"
$case.Range("E28").Value = "hold_correction/hold_correction3"

$case.Range("A26").Value = "24"
$case.Range("A27").Value = "25"
$case.Range("A28").Value = "26"

# New "mpar" marker column for the rows that previously lacked it
$case.Range("Q19:Q25").Value = "mpar"

$case.Activate()
$case.Range("P26").Select()
$suite.Activate()
$suite.Range("B14").Select()
